$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nbsp = [char]0x00A0
$needle = "," + $nbsp + ","
$replacement = "," + $nbsp + "Copiar el horario,"

for ($row = 2; $row -le 6; $row++) {
    # Remove the stray ", ," (comma + nbsp + comma) left after the first
    # day's schedule in the structuredData HTML, replacing it with the
    # missing "Copiar el horario" label (column M).
    $mCell = $ws.Cells.Item($row, 13)
    $mValue = $mCell.Value()
    if ($mValue -ne $null -and $mValue.Contains($needle)) {
        $mCell.Value = $mValue.Replace($needle, $replacement)
    }

    # Drop the now-unused "horario" column F cell entirely for this row.
    $ws.Cells.Item($row, 6).Value = $null
}
